# Auto-generated Excel COM-interop script applying the Chocobo_Profits diff.
# Updates currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfit columns
# (H..N) for the affected leve rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")  # row 70
$ws.Cells.Item(70, 8).Value = 2325.5881  # H70: 2104.75 -> 2325.5881
$ws.Cells.Item(70, 9).Value = 887.2727  # I70: 886.1539 -> 887.2727
$ws.Cells.Item(70, 10).Value = 4962.5  # J70: 4367.857 -> 4962.5
$ws.Cells.Item(70, 11).Value = 2661.8181  # K70: 2658.4617 -> 2661.8181
$ws.Cells.Item(70, 12).Value = 14887.5  # L70: 13103.571 -> 14887.5
$ws.Cells.Item(70, 13).Value = -2391.8181  # M70: -2388.4617 -> -2391.8181
$ws.Cells.Item(70, 14).Value = -15427.5  # N70: -13643.571 -> -15427.5

$ws = $wb.Worksheets.Item("ALC")  # row 73
$ws.Cells.Item(73, 8).Value = 2325.5881  # H73: 2104.75 -> 2325.5881
$ws.Cells.Item(73, 9).Value = 887.2727  # I73: 886.1539 -> 887.2727
$ws.Cells.Item(73, 10).Value = 4962.5  # J73: 4367.857 -> 4962.5
$ws.Cells.Item(73, 11).Value = 2661.8181  # K73: 2658.4617 -> 2661.8181
$ws.Cells.Item(73, 12).Value = 14887.5  # L73: 13103.571 -> 14887.5
$ws.Cells.Item(73, 13).Value = -1725.8181  # M73: -1722.4617 -> -1725.8181
$ws.Cells.Item(73, 14).Value = -16759.5  # N73: -14975.571 -> -16759.5

$ws = $wb.Worksheets.Item("ALC")  # row 109
$ws.Cells.Item(109, 8).Value = 27190.477  # H109: 30466.666 -> 27190.477
$ws.Cells.Item(109, 10).Value = 27190.477  # J109: 30466.666 -> 27190.477
$ws.Cells.Item(109, 12).Value = 27190.477  # L109: 30466.666 -> 27190.477
$ws.Cells.Item(109, 14).Value = -29964.477  # N109: -33240.666 -> -29964.477

$ws = $wb.Worksheets.Item("ALC")  # row 113
$ws.Cells.Item(113, 8).Value = 8113.1816  # H113: 11842.5 -> 8113.1816
$ws.Cells.Item(113, 9).Value = 4964.8887  # I113: 6501.3335 -> 4964.8887
$ws.Cells.Item(113, 10).Value = 10292.77  # J113: 13622.889 -> 10292.77
$ws.Cells.Item(113, 11).Value = 4964.8887  # K113: 6501.3335 -> 4964.8887
$ws.Cells.Item(113, 12).Value = 10292.77  # L113: 13622.889 -> 10292.77
$ws.Cells.Item(113, 13).Value = -1710.8887  # M113: -3247.3335 -> -1710.8887
$ws.Cells.Item(113, 14).Value = -16800.77  # N113: -20130.889 -> -16800.77

$ws = $wb.Worksheets.Item("ALC")  # row 116
$ws.Cells.Item(116, 8).Value = 8000.7144  # H116: 8300.75 -> 8000.7144
$ws.Cells.Item(116, 9).Value = 2319.5  # I116: 2479.25 -> 2319.5
$ws.Cells.Item(116, 10).Value = 10273.2  # J116: 9756.125 -> 10273.2
$ws.Cells.Item(116, 11).Value = 2319.5  # K116: 2479.25 -> 2319.5
$ws.Cells.Item(116, 12).Value = 10273.2  # L116: 9756.125 -> 10273.2
$ws.Cells.Item(116, 13).Value = 1122.5  # M116: 962.75 -> 1122.5
$ws.Cells.Item(116, 14).Value = -17157.2  # N116: -16640.125 -> -17157.2

$ws = $wb.Worksheets.Item("ALC")  # row 132
$ws.Cells.Item(132, 8).Value = 6373.3  # H132: 5538.0586 -> 6373.3
$ws.Cells.Item(132, 9).Value = 5855.72  # I132: 4949.6 -> 5855.72
$ws.Cells.Item(132, 10).Value = 8961.200000000001  # J132: 9951.5 -> 8961.200000000001
$ws.Cells.Item(132, 11).Value = 17567.16  # K132: 14848.8 -> 17567.16
$ws.Cells.Item(132, 12).Value = 26883.6  # L132: 29854.5 -> 26883.6
$ws.Cells.Item(132, 13).Value = -15037.16  # M132: -12318.8 -> -15037.16
$ws.Cells.Item(132, 14).Value = -31943.6  # N132: -34914.5 -> -31943.6

$ws = $wb.Worksheets.Item("ARM")  # row 32
$ws.Cells.Item(32, 8).Value = 6129.093  # H32: 6127.4443 -> 6129.093
$ws.Cells.Item(32, 9).Value = 3555.2058  # I32: 3366.3333 -> 3555.2058
$ws.Cells.Item(32, 10).Value = 10504.7  # J32: 10466.333 -> 10504.7
$ws.Cells.Item(32, 11).Value = 3555.2058  # K32: 3366.3333 -> 3555.2058
$ws.Cells.Item(32, 12).Value = 10504.7  # L32: 10466.333 -> 10504.7
$ws.Cells.Item(32, 13).Value = -3268.2058  # M32: -3079.3333 -> -3268.2058
$ws.Cells.Item(32, 14).Value = -11078.7  # N32: -11040.333 -> -11078.7

$ws = $wb.Worksheets.Item("ARM")  # row 43
$ws.Cells.Item(43, 8).Value = 21715.2  # H43: 21388.5 -> 21715.2
$ws.Cells.Item(43, 10).Value = 24644  # J43: 32777 -> 24644
$ws.Cells.Item(43, 12).Value = 24644  # L43: 32777 -> 24644
$ws.Cells.Item(43, 14).Value = -25270  # N43: -33403 -> -25270

$ws = $wb.Worksheets.Item("ARM")  # row 74
$ws.Cells.Item(74, 8).Value = 1687.075  # H74: 1752.3846 -> 1687.075
$ws.Cells.Item(74, 9).Value = 1300.9062  # I74: 1387 -> 1300.9062
$ws.Cells.Item(74, 10).Value = 3231.75  # J74: 2970.3333 -> 3231.75
$ws.Cells.Item(74, 11).Value = 1300.9062  # K74: 1387 -> 1300.9062
$ws.Cells.Item(74, 12).Value = 3231.75  # L74: 2970.3333 -> 3231.75
$ws.Cells.Item(74, 13).Value = -426.9061999999999  # M74: -513 -> -426.9061999999999
$ws.Cells.Item(74, 14).Value = -4979.75  # N74: -4718.3333 -> -4979.75

$ws = $wb.Worksheets.Item("ARM")  # row 77
$ws.Cells.Item(77, 8).Value = 1687.075  # H77: 1752.3846 -> 1687.075
$ws.Cells.Item(77, 9).Value = 1300.9062  # I77: 1387 -> 1300.9062
$ws.Cells.Item(77, 10).Value = 3231.75  # J77: 2970.3333 -> 3231.75
$ws.Cells.Item(77, 11).Value = 6504.530999999999  # K77: 6935 -> 6504.530999999999
$ws.Cells.Item(77, 12).Value = 16158.75  # L77: 14851.6665 -> 16158.75
$ws.Cells.Item(77, 13).Value = -2136.530999999999  # M77: -2567 -> -2136.530999999999
$ws.Cells.Item(77, 14).Value = -24894.75  # N77: -23587.6665 -> -24894.75

$ws = $wb.Worksheets.Item("ARM")  # row 122
$ws.Cells.Item(122, 8).Value = 2448.1304  # H122: 2904.375 -> 2448.1304
$ws.Cells.Item(122, 9).Value = 1621.9333  # I122: 1771.909 -> 1621.9333
$ws.Cells.Item(122, 10).Value = 3997.25  # J122: 5395.8 -> 3997.25
$ws.Cells.Item(122, 11).Value = 4865.7999  # K122: 5315.727000000001 -> 4865.7999
$ws.Cells.Item(122, 12).Value = 11991.75  # L122: 16187.4 -> 11991.75
$ws.Cells.Item(122, 13).Value = -2415.7999  # M122: -2865.727000000001 -> -2415.7999
$ws.Cells.Item(122, 14).Value = -16891.75  # N122: -21087.4 -> -16891.75

$ws = $wb.Worksheets.Item("ARM")  # row 132
$ws.Cells.Item(132, 8).Value = 2242  # H132: 2295.9792 -> 2242
$ws.Cells.Item(132, 9).Value = 1774.4857  # I132: 1771.0278 -> 1774.4857
$ws.Cells.Item(132, 10).Value = 3410.7856  # J132: 3870.8333 -> 3410.7856
$ws.Cells.Item(132, 11).Value = 5323.4571  # K132: 5313.0834 -> 5323.4571
$ws.Cells.Item(132, 12).Value = 10232.3568  # L132: 11612.4999 -> 10232.3568
$ws.Cells.Item(132, 13).Value = -2793.4571  # M132: -2783.0834 -> -2793.4571
$ws.Cells.Item(132, 14).Value = -15292.3568  # N132: -16672.4999 -> -15292.3568

$ws = $wb.Worksheets.Item("BSM")  # row 134
$ws.Cells.Item(134, 8).Value = 2475.94  # H134: 2724.2974 -> 2475.94
$ws.Cells.Item(134, 9).Value = 1390.8422  # I134: 1229 -> 1390.8422
$ws.Cells.Item(134, 10).Value = 5912.0835  # J134: 6258.636 -> 5912.0835
$ws.Cells.Item(134, 11).Value = 4172.5266  # K134: 3687 -> 4172.5266
$ws.Cells.Item(134, 12).Value = 17736.2505  # L134: 18775.908 -> 17736.2505
$ws.Cells.Item(134, 13).Value = -1637.5266  # M134: -1152 -> -1637.5266
$ws.Cells.Item(134, 14).Value = -22806.2505  # N134: -23845.908 -> -22806.2505

$ws = $wb.Worksheets.Item("CRP")  # row 58
$ws.Cells.Item(58, 8).Value = 2050  # H58: 2371.32 -> 2050
$ws.Cells.Item(58, 9).Value = 1678.7678  # I58: 1244.6923 -> 1678.7678
$ws.Cells.Item(58, 10).Value = 4128.9  # J58: 3591.8333 -> 4128.9
$ws.Cells.Item(58, 11).Value = 1678.7678  # K58: 1244.6923 -> 1678.7678
$ws.Cells.Item(58, 12).Value = 4128.9  # L58: 3591.8333 -> 4128.9
$ws.Cells.Item(58, 13).Value = -1475.7678  # M58: -1041.6923 -> -1475.7678
$ws.Cells.Item(58, 14).Value = -4534.9  # N58: -3997.8333 -> -4534.9

$ws = $wb.Worksheets.Item("CRP")  # row 136
$ws.Cells.Item(136, 8).Value = 2050  # H136: 2371.32 -> 2050
$ws.Cells.Item(136, 9).Value = 1678.7678  # I136: 1244.6923 -> 1678.7678
$ws.Cells.Item(136, 10).Value = 4128.9  # J136: 3591.8333 -> 4128.9
$ws.Cells.Item(136, 11).Value = 5036.303400000001  # K136: 3734.0769 -> 5036.303400000001
$ws.Cells.Item(136, 12).Value = 12386.7  # L136: 10775.4999 -> 12386.7
$ws.Cells.Item(136, 13).Value = -2486.303400000001  # M136: -1184.0769 -> -2486.303400000001
$ws.Cells.Item(136, 14).Value = -17486.7  # N136: -15875.4999 -> -17486.7

$ws = $wb.Worksheets.Item("CRP")  # row 141
$ws.Cells.Item(141, 8).Value = 24511.111  # H141: 24763.637 -> 24511.111
$ws.Cells.Item(141, 10).Value = 24511.111  # J141: 24763.637 -> 24511.111
$ws.Cells.Item(141, 12).Value = 24511.111  # L141: 24763.637 -> 24511.111
$ws.Cells.Item(141, 14).Value = -34871.111  # N141: -35123.637 -> -34871.111

$ws = $wb.Worksheets.Item("CUL")  # row 74
$ws.Cells.Item(74, 8).Value = 10498.5  # H74: 9669 -> 10498.5
$ws.Cells.Item(74, 9).Value = 8000  # I74: 5013 -> 8000
$ws.Cells.Item(74, 10).Value = 11331.333  # J74: 11997 -> 11331.333
$ws.Cells.Item(74, 11).Value = 24000  # K74: 15039 -> 24000
$ws.Cells.Item(74, 12).Value = 33993.999  # L74: 35991 -> 33993.999
$ws.Cells.Item(74, 13).Value = -22939  # M74: -13978 -> -22939
$ws.Cells.Item(74, 14).Value = -36115.999  # N74: -38113 -> -36115.999

$ws = $wb.Worksheets.Item("CUL")  # row 75
$ws.Cells.Item(75, 8).Value = 2930.4285  # H75: 2315.6667 -> 2930.4285
$ws.Cells.Item(75, 9).Value = 1013  # I75: 914 -> 1013
$ws.Cells.Item(75, 10).Value = 3250  # J75: 2596 -> 3250
$ws.Cells.Item(75, 11).Value = 3039  # K75: 2742 -> 3039
$ws.Cells.Item(75, 12).Value = 9750  # L75: 7788 -> 9750
$ws.Cells.Item(75, 13).Value = -2041  # M75: -1744 -> -2041
$ws.Cells.Item(75, 14).Value = -11746  # N75: -9784 -> -11746

$ws = $wb.Worksheets.Item("CUL")  # row 76
$ws.Cells.Item(76, 8).Value = 3714.2856  # H76: 3600 -> 3714.2856
$ws.Cells.Item(76, 10).Value = 5000  # J76: 4666.6665 -> 5000
$ws.Cells.Item(76, 12).Value = 15000  # L76: 13999.9995 -> 15000
$ws.Cells.Item(76, 14).Value = -15766  # N76: -14765.9995 -> -15766

$ws = $wb.Worksheets.Item("CUL")  # row 77
$ws.Cells.Item(77, 8).Value = 10498.5  # H77: 9669 -> 10498.5
$ws.Cells.Item(77, 9).Value = 8000  # I77: 5013 -> 8000
$ws.Cells.Item(77, 10).Value = 11331.333  # J77: 11997 -> 11331.333
$ws.Cells.Item(77, 11).Value = 72000  # K77: 45117 -> 72000
$ws.Cells.Item(77, 12).Value = 101981.997  # L77: 107973 -> 101981.997
$ws.Cells.Item(77, 13).Value = -66696  # M77: -39813 -> -66696
$ws.Cells.Item(77, 14).Value = -112589.997  # N77: -118581 -> -112589.997

$ws = $wb.Worksheets.Item("CUL")  # row 78
$ws.Cells.Item(78, 8).Value = 2930.4285  # H78: 2315.6667 -> 2930.4285
$ws.Cells.Item(78, 9).Value = 1013  # I78: 914 -> 1013
$ws.Cells.Item(78, 10).Value = 3250  # J78: 2596 -> 3250
$ws.Cells.Item(78, 11).Value = 9117  # K78: 8226 -> 9117
$ws.Cells.Item(78, 12).Value = 29250  # L78: 23364 -> 29250
$ws.Cells.Item(78, 13).Value = -4125  # M78: -3234 -> -4125
$ws.Cells.Item(78, 14).Value = -39234  # N78: -33348 -> -39234

$ws = $wb.Worksheets.Item("CUL")  # row 79
$ws.Cells.Item(79, 8).Value = 3714.2856  # H79: 3600 -> 3714.2856
$ws.Cells.Item(79, 10).Value = 5000  # J79: 4666.6665 -> 5000
$ws.Cells.Item(79, 12).Value = 15000  # L79: 13999.9995 -> 15000
$ws.Cells.Item(79, 14).Value = -17652  # N79: -16651.9995 -> -17652

$ws = $wb.Worksheets.Item("CUL")  # row 80
$ws.Cells.Item(80, 8).Value = 9599.4  # H80: 9799.200000000001 -> 9599.4
$ws.Cells.Item(80, 10).Value = 9599.4  # J80: 9799.200000000001 -> 9599.4
$ws.Cells.Item(80, 12).Value = 28798.2  # L80: 29397.6 -> 28798.2
$ws.Cells.Item(80, 14).Value = -30670.2  # N80: -31269.6 -> -30670.2

$ws = $wb.Worksheets.Item("CUL")  # row 82
$ws.Cells.Item(82, 8).Value = 3834.7778  # H82: 6218.8335 -> 3834.7778
$ws.Cells.Item(82, 9).Value = 682.6  # I82: 756.5 -> 682.6
$ws.Cells.Item(82, 10).Value = 7775  # J82: 8950 -> 7775
$ws.Cells.Item(82, 11).Value = 2047.8  # K82: 2269.5 -> 2047.8
$ws.Cells.Item(82, 12).Value = 23325  # L82: 26850 -> 23325
$ws.Cells.Item(82, 13).Value = -1641.8  # M82: -1863.5 -> -1641.8
$ws.Cells.Item(82, 14).Value = -24137  # N82: -27662 -> -24137

$ws = $wb.Worksheets.Item("CUL")  # row 83
$ws.Cells.Item(83, 8).Value = 9599.4  # H83: 9799.200000000001 -> 9599.4
$ws.Cells.Item(83, 10).Value = 9599.4  # J83: 9799.200000000001 -> 9599.4
$ws.Cells.Item(83, 12).Value = 86394.59999999999  # L83: 88192.8 -> 86394.59999999999
$ws.Cells.Item(83, 14).Value = -95754.59999999999  # N83: -97552.8 -> -95754.59999999999

$ws = $wb.Worksheets.Item("CUL")  # row 85
$ws.Cells.Item(85, 8).Value = 3834.7778  # H85: 6218.8335 -> 3834.7778
$ws.Cells.Item(85, 9).Value = 682.6  # I85: 756.5 -> 682.6
$ws.Cells.Item(85, 10).Value = 7775  # J85: 8950 -> 7775
$ws.Cells.Item(85, 11).Value = 2047.8  # K85: 2269.5 -> 2047.8
$ws.Cells.Item(85, 12).Value = 23325  # L85: 26850 -> 23325
$ws.Cells.Item(85, 13).Value = -643.8000000000002  # M85: -865.5 -> -643.8000000000002
$ws.Cells.Item(85, 14).Value = -26133  # N85: -29658 -> -26133

$ws = $wb.Worksheets.Item("CUL")  # row 86
$ws.Cells.Item(86, 8).Value = 30127  # H86: 30252.25 -> 30127
$ws.Cells.Item(86, 9).Value = 502  # I86: 0 -> 502
$ws.Cells.Item(86, 10).Value = 40002  # J86: 30252.25 -> 40002
$ws.Cells.Item(86, 11).Value = 1506  # K86: 0 -> 1506
$ws.Cells.Item(86, 12).Value = 120006  # L86: 90756.75 -> 120006
$ws.Cells.Item(86, 13).Value = -320  # M86: None -> -320
$ws.Cells.Item(86, 14).Value = -122378  # N86: -93128.75 -> -122378

$ws = $wb.Worksheets.Item("CUL")  # row 87
$ws.Cells.Item(87, 8).Value = 750  # H87: 3750 -> 750
$ws.Cells.Item(87, 9).Value = 750  # I87: 3750 -> 750
$ws.Cells.Item(87, 11).Value = 2250  # K87: 11250 -> 2250
$ws.Cells.Item(87, 13).Value = -1002  # M87: -10002 -> -1002

$ws = $wb.Worksheets.Item("CUL")  # row 89
$ws.Cells.Item(89, 8).Value = 30127  # H89: 30252.25 -> 30127
$ws.Cells.Item(89, 9).Value = 502  # I89: 0 -> 502
$ws.Cells.Item(89, 10).Value = 40002  # J89: 30252.25 -> 40002
$ws.Cells.Item(89, 11).Value = 4518  # K89: 0 -> 4518
$ws.Cells.Item(89, 12).Value = 360018  # L89: 272270.25 -> 360018
$ws.Cells.Item(89, 13).Value = 1410  # M89: None -> 1410
$ws.Cells.Item(89, 14).Value = -371874  # N89: -284126.25 -> -371874

$ws = $wb.Worksheets.Item("CUL")  # row 90
$ws.Cells.Item(90, 8).Value = 750  # H90: 3750 -> 750
$ws.Cells.Item(90, 9).Value = 750  # I90: 3750 -> 750
$ws.Cells.Item(90, 11).Value = 6750  # K90: 33750 -> 6750
$ws.Cells.Item(90, 13).Value = -510  # M90: -27510 -> -510

$ws = $wb.Worksheets.Item("GSM")  # row 26
$ws.Cells.Item(26, 8).Value = 20000  # H26: 12440 -> 20000
$ws.Cells.Item(26, 9).Value = 0  # I26: 13880 -> 0
$ws.Cells.Item(26, 10).Value = 20000  # J26: 11000 -> 20000
$ws.Cells.Item(26, 11).Value = 0  # K26: 13880 -> 0
$ws.Cells.Item(26, 12).Value = 20000  # L26: 11000 -> 20000
$ws.Cells.Item(26, 13).ClearContents()  # M26: -13600 -> (removed)
$ws.Cells.Item(26, 14).Value = -20560  # N26: -11560 -> -20560

$ws = $wb.Worksheets.Item("GSM")  # row 48
$ws.Cells.Item(48, 8).Value = 35000  # H48: 0 -> 35000
$ws.Cells.Item(48, 10).Value = 35000  # J48: 0 -> 35000
$ws.Cells.Item(48, 12).Value = 35000  # L48: 0 -> 35000
$ws.Cells.Item(48, 13).Value = -35970  # M48: None -> -35970

$ws = $wb.Worksheets.Item("GSM")  # row 50
$ws.Cells.Item(50, 8).Value = 20000  # H50: 12440 -> 20000
$ws.Cells.Item(50, 9).Value = 0  # I50: 13880 -> 0
$ws.Cells.Item(50, 10).Value = 20000  # J50: 11000 -> 20000
$ws.Cells.Item(50, 11).Value = 0  # K50: 13880 -> 0
$ws.Cells.Item(50, 12).Value = 20000  # L50: 11000 -> 20000
$ws.Cells.Item(50, 13).ClearContents()  # M50: -13382 -> (removed)
$ws.Cells.Item(50, 14).Value = -20996  # N50: -11996 -> -20996

$ws = $wb.Worksheets.Item("LTW")  # row 122
$ws.Cells.Item(122, 8).Value = 4651.92  # H122: 3883.303 -> 4651.92
$ws.Cells.Item(122, 9).Value = 2954.5334  # I122: 2513.45 -> 2954.5334
$ws.Cells.Item(122, 10).Value = 7198  # J122: 5990.769 -> 7198
$ws.Cells.Item(122, 11).Value = 8863.600199999999  # K122: 7540.349999999999 -> 8863.600199999999
$ws.Cells.Item(122, 12).Value = 21594  # L122: 17972.307 -> 21594
$ws.Cells.Item(122, 13).Value = -6413.600199999999  # M122: -5090.349999999999 -> -6413.600199999999
$ws.Cells.Item(122, 14).Value = -26494  # N122: -22872.307 -> -26494

$ws = $wb.Worksheets.Item("LTW")  # row 132
$ws.Cells.Item(132, 8).Value = 2975.6902  # H132: 2953.4226 -> 2975.6902
$ws.Cells.Item(132, 9).Value = 1227.9565  # I132: 1253.1915 -> 1227.9565
$ws.Cells.Item(132, 10).Value = 6191.52  # J132: 6283.0415 -> 6191.52
$ws.Cells.Item(132, 11).Value = 3683.8695  # K132: 3759.5745 -> 3683.8695
$ws.Cells.Item(132, 12).Value = 18574.56  # L132: 18849.1245 -> 18574.56
$ws.Cells.Item(132, 13).Value = -1153.8695  # M132: -1229.5745 -> -1153.8695
$ws.Cells.Item(132, 14).Value = -23634.56  # N132: -23909.1245 -> -23634.56
